# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (header in G1) previously held "Strike#" values; this
# script rewrites the column G data (rows 2-15) with the recalculated
# strikeout ("K") values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(4, 4, 5, 6, 7, 5, 2, 7, 2, 4, 2, 5, 2, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
